# Add the new 2021 data row (row 10) to the worksheet, copying the
# formatting of the existing last data row (row 9) so the new row
# matches the look of the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 9's formatting (bold, bordered, centered style on col A;
# plain style on col B) down onto row 10.
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 1158
